$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("PDiBCpDoC")

# --- Remove the embedded chart picture on the About sheet ---
while ($wsAbout.Shapes.Count -gt 0) {
    $wsAbout.Shapes.Item(1).Delete()
}

# --- About sheet text / source citation updates ---
$wsAbout.Range("B3").Value = "Massachusetts Institute of Technology"
$wsAbout.Range("B4").Value = 2021
$wsAbout.Range("B5").Value = "Re-examining rates of lithium-ion battery technology improvement and cost decline"
$wsAbout.Range("B6").Value = "https://pubs.rsc.org/en/content/articlepdf/2021/ee/d0ee02681f?page=search"
$wsAbout.Range("B7").Value = "Abstract"
$wsAbout.Range("C8").Value = ""
$wsAbout.Range("A9").Value = "Note: We take the average of learning rates quoted in the Abstract (20%-27%)"

# --- PDiBCpDoC sheet: label + recomputed value ---
$wsData.Range("B1").Value = "Perc Decline per Doubling (dimensionless)"
$wsData.Range("B2").Formula = "=AVERAGE(0.2,0.27)"
